# Update "想去人数" (interested-count) figures in column F on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 295
$ws1.Range("F3").Value = 203
$ws1.Range("F4").Value = 2466
$ws1.Range("F5").Value = 1802
$ws1.Range("F6").Value = 342
$ws1.Range("F7").Value = 100
$ws1.Range("F8").Value = 852
$ws1.Range("F9").Value = 170

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 295
$ws4.Range("F3").Value = 203
$ws4.Range("F4").Value = 2466
$ws4.Range("F5").Value = 1802
$ws4.Range("F6").Value = 342
$ws4.Range("F8").Value = 100
$ws4.Range("F9").Value = 852
$ws4.Range("F10").Value = 170
